$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.475.00'
$ws.Range('E2').Value = '  -3.62%  '
$ws.Range('D3').Value = '2.593.16'
$ws.Range('E3').Value = '  -2.70%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.03'
$ws.Range('E5').Value = '  -4.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.85'
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -4.64%  '
$ws.Range('E10').Value = '  -0.42%  '
$ws.Range('E11').Value = '  -5.79%  '
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.09'
$ws.Range('E13').Value = '  -3.39%  '
$ws.Range('D14').Value = '3.058.53'
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('E15').Value = '  -8.95%  '
$ws.Range('D16').Value = '63.290.83'
$ws.Range('E16').Value = '  -3.69%  '
$ws.Range('D17').Value = '2.585.17'
$ws.Range('E17').Value = '  -3.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.98'
$ws.Range('E18').Value = '  -5.06%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.54'
$ws.Range('E19').Value = '  -5.69%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.48'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '341.60'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.32'
$ws.Range('E23').Value = '  -3.69%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  -4.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.12'
$ws.Range('E26').Value = '  -5.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '575.93'
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.56'
$ws.Range('E28').Value = '  -3.92%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.161'
$ws.Range('E30').Value = '  -1.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.87'
$ws.Range('E31').Value = '  -3.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.06'
$ws.Range('E32').Value = '  -4.66%  '
$ws.Range('E33').Value = '  -5.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.53'
$ws.Range('E34').Value = '  -3.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.43'
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('E36').Value = '  -5.17%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  -4.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '154.73'
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.86'
$ws.Range('E40').Value = '  -5.45%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.49'
$ws.Range('E42').Value = '  +7.21%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.26'
$ws.Range('E43').Value = '  -3.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '156.53'
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('E45').Value = '  -5.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.15'
$ws.Range('E46').Value = '  -0.25%  '
$ws.Range('E47').Value = '  -5.41%  '
$ws.Range('E48').Value = '  -3.00%  '
$ws.Range('E49').Value = '  -2.39%  '
$ws.Range('E50').Value = '  -5.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.73'
$ws.Range('E51').Value = '  -5.45%  '
